$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "pontos" column (D) holds numeric-looking values that are stored as
# text in the source data, so force a text number format before writing
# them to keep them as text (not auto-converted to numbers).
$ws.Range("D2:D9").NumberFormat = "@"

# Update existing rows (C: nick name, D: pontos)
$ws.Range("C2").Value = "bona final boss"
$ws.Range("D2").Value = "9999"

$ws.Range("C3").Value = "bona semi boss"
$ws.Range("D3").Value = "2501"

$ws.Range("C4").Value = "diego"
$ws.Range("D4").Value = "1935"

$ws.Range("C5").Value = "diego"
$ws.Range("D5").Value = "1891"

$ws.Range("C6").Value = "diego"
$ws.Range("D6").Value = "1826"

$ws.Range("C7").Value = "diego"
$ws.Range("D7").Value = "1791"

# Add new row 8, copying style from A7 into A8 (keeps the bordered/bold style)
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "7º"
$ws.Range("C8").Value = "diego"
$ws.Range("D8").Value = "1424"

# Add new row 9, copying style from A7 into A9
$ws.Range("A7").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "8º"
$ws.Range("C9").Value = "diego"
$ws.Range("D9").Value = "1423"
